$d = $word.ActiveDocument

# 1) Номер аудитории: 44 -> 43
$d.Content.Find.Execute("44", $true, $false, $false, $false, $false,
                         $true, 1, $false, "43", 2)

# 2) Отчетный период: с 24.04.2024 по 08.05.2024 -> с 18.05.2024 по 21.05.2024
$d.Content.Find.Execute("с 24.04.2024 по 08.05.2024", $true, $false, $false, $false, $false,
                         $true, 1, $false, "с 18.05.2024 по 21.05.2024", 2)

# 3) Table cell: ВСЕ ПК -> sdfg
$d.Content.Find.Execute("ВСЕ ПК", $true, $false, $false, $false, $false,
                         $true, 1, $false, "sdfg", 2)

# 4) Table cell: не работают мышки и кнопки включения -> dsfgdsfg
$d.Content.Find.Execute("не работают мышки и кнопки включения", $true, $false, $false, $false, $false,
                         $true, 1, $false, "dsfgdsfg", 2)
